$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 308.6
$ws.Range("I28").Value = 308.6
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 308.6
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 176.4
$ws.Range("N28").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 62506040
$ws.Range("I62").Value = 90911510
$ws.Range("J62").Value = 14000
$ws.Range("K62").Value = 90911510
$ws.Range("L62").Value = 14000
$ws.Range("M62").Value = -90910886
$ws.Range("N62").Value = -15248

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 62506040
$ws.Range("I65").Value = 90911510
$ws.Range("J65").Value = 14000
$ws.Range("K65").Value = 454557550
$ws.Range("L65").Value = 70000
$ws.Range("M65").Value = -454554430
$ws.Range("N65").Value = -76240

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 7937258.5
$ws.Range("I92").Value = 13889224
$ws.Range("J92").Value = 1305
$ws.Range("K92").Value = 13889224
$ws.Range("L92").Value = 1305
$ws.Range("M92").Value = -13887976
$ws.Range("N92").Value = -3801

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 5142.857
$ws.Range("I106").Value = 5333.3335
$ws.Range("J106").Value = 4000
$ws.Range("K106").Value = 5333.3335
$ws.Range("L106").Value = 4000
$ws.Range("M106").Value = -4702.3335
$ws.Range("N106").Value = -5262

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 173.9375
$ws.Range("I107").Value = 165.53334
$ws.Range("J107").Value = 300
$ws.Range("K107").Value = 165.53334
$ws.Range("L107").Value = 300
$ws.Range("M107").Value = 1754.46666
$ws.Range("N107").Value = -4140

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4239.0356
$ws.Range("I132").Value = 4307.72
$ws.Range("J132").Value = 3666.6667
$ws.Range("K132").Value = 12923.16
$ws.Range("L132").Value = 11000.0001
$ws.Range("M132").Value = -10393.16
$ws.Range("N132").Value = -16060.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2573.42
$ws.Range("I32").Value = 1393.7303
$ws.Range("K32").Value = 1393.7303
$ws.Range("M32").Value = -1106.7303

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1749.4642
$ws.Range("I45").Value = 1792.5
$ws.Range("J45").Value = 1641.875
$ws.Range("K45").Value = 1792.5
$ws.Range("L45").Value = 1641.875
$ws.Range("M45").Value = -1415.5
$ws.Range("N45").Value = -2395.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1191.683
$ws.Range("I61").Value = 1074.069
$ws.Range("J61").Value = 1475.9166
$ws.Range("K61").Value = 1074.069
$ws.Range("L61").Value = 1475.9166
$ws.Range("M61").Value = -862.069
$ws.Range("N61").Value = -1899.9166

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1191.683
$ws.Range("I136").Value = 1074.069
$ws.Range("J136").Value = 1475.9166
$ws.Range("K136").Value = 3222.207
$ws.Range("L136").Value = 4427.7498
$ws.Range("M136").Value = -672.2069999999999
$ws.Range("N136").Value = -9527.7498

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3237.8076
$ws.Range("I20").Value = 1510.9166
$ws.Range("J20").Value = 4718
$ws.Range("K20").Value = 1510.9166
$ws.Range("L20").Value = 4718
$ws.Range("M20").Value = -1263.9166
$ws.Range("N20").Value = -5212

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1943.5555
$ws.Range("I107").Value = 1581
$ws.Range("J107").Value = 2513.2856
$ws.Range("K107").Value = 1581
$ws.Range("L107").Value = 2513.2856
$ws.Range("M107").Value = 339
$ws.Range("N107").Value = -6353.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3655829.5
$ws.Range("I31").Value = 6282204
$ws.Range("J31").Value = 1743.4348
$ws.Range("K31").Value = 6282204
$ws.Range("L31").Value = 1743.4348
$ws.Range("M31").Value = -6281909
$ws.Range("N31").Value = -2333.4348

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3655829.5
$ws.Range("I34").Value = 6282204
$ws.Range("J34").Value = 1743.4348
$ws.Range("K34").Value = 6282204
$ws.Range("L34").Value = 1743.4348
$ws.Range("M34").Value = -6282002
$ws.Range("N34").Value = -2147.4348

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3974473.8
$ws.Range("I99").Value = 8937666
$ws.Range("J99").Value = 3920
$ws.Range("K99").Value = 8937666
$ws.Range("L99").Value = 3920
$ws.Range("M99").Value = -8936168
$ws.Range("N99").Value = -6916

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3974473.8
$ws.Range("I126").Value = 8937666
$ws.Range("J126").Value = 3920
$ws.Range("K126").Value = 26812998
$ws.Range("L126").Value = 11760
$ws.Range("M126").Value = -26810528
$ws.Range("N126").Value = -16700

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3446.1538
$ws.Range("I134").Value = 3640
$ws.Range("J134").Value = 2800
$ws.Range("K134").Value = 10920
$ws.Range("L134").Value = 8400
$ws.Range("M134").Value = -8385
$ws.Range("N134").Value = -13470

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 836.3889
$ws.Range("I131").Value = 461
$ws.Range("J131").Value = 980.7692
$ws.Range("K131").Value = 1383
$ws.Range("L131").Value = 2942.3076
$ws.Range("M131").Value = 3657
$ws.Range("N131").Value = -13022.3076

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 314.16
$ws.Range("I107").Value = 128.23077
$ws.Range("J107").Value = 515.5833
$ws.Range("K107").Value = 128.23077
$ws.Range("L107").Value = 515.5833
$ws.Range("M107").Value = 1791.76923
$ws.Range("N107").Value = -4355.5833

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2636.92
$ws.Range("I132").Value = 2216.4546
$ws.Range("J132").Value = 2967.2856
$ws.Range("K132").Value = 6649.3638
$ws.Range("L132").Value = 8901.856800000001
$ws.Range("M132").Value = -4119.3638
$ws.Range("N132").Value = -13961.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 773.0625
$ws.Range("I22").Value = 491.26666
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 491.26666
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = -196.26666
$ws.Range("N22").Value = -5590

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 773.0625
$ws.Range("I27").Value = 491.26666
$ws.Range("J27").Value = 5000
$ws.Range("K27").Value = 491.26666
$ws.Range("L27").Value = 5000
$ws.Range("M27").Value = -384.26666
$ws.Range("N27").Value = -5214

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 679.9286
$ws.Range("I46").Value = 618.8889
$ws.Range("J46").Value = 789.8
$ws.Range("K46").Value = 618.8889
$ws.Range("L46").Value = 789.8
$ws.Range("M46").Value = -430.8889
$ws.Range("N46").Value = -1165.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1981.0952
$ws.Range("I61").Value = 1909
$ws.Range("J61").Value = 2077.2222
$ws.Range("K61").Value = 1909
$ws.Range("L61").Value = 2077.2222
$ws.Range("M61").Value = -1707
$ws.Range("N61").Value = -2481.2222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1981.0952
$ws.Range("I113").Value = 1909
$ws.Range("J113").Value = 2077.2222
$ws.Range("K113").Value = 1909
$ws.Range("L113").Value = 2077.2222
$ws.Range("M113").Value = 261
$ws.Range("N113").Value = -6417.2222

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 660.08
$ws.Range("I113").Value = 642.9048
$ws.Range("J113").Value = 750.25
$ws.Range("K113").Value = 1928.7144
$ws.Range("L113").Value = 2250.75
$ws.Range("M113").Value = 241.2855999999999
$ws.Range("N113").Value = -6590.75
